$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 409-410, pushing the existing data
# (old rows 409..484) down to 411..486.
$ws.Range("A409:A410").EntireRow.Insert()

# New row 409: "Primera" quality entry for date 44785 (2022-08-12)
$ws.Cells.Item(409, 1).Value  = 8
$ws.Cells.Item(409, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(409, 3).Value  = "Coquimbo"
$ws.Cells.Item(409, 4).Value  = 44785
$ws.Cells.Item(409, 5).Value  = 4
$ws.Cells.Item(409, 6).Value  = 100112043
$ws.Cells.Item(409, 7).Value  = "Pepino dulce"
$ws.Cells.Item(409, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(409, 9).Value  = "Primera"
$ws.Cells.Item(409, 10).Value = 520
$ws.Cells.Item(409, 11).Value = 13000
$ws.Cells.Item(409, 12).Value = 13500
$ws.Cells.Item(409, 13).Value = 13250
$ws.Cells.Item(409, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(409, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(409, 16).Value = 736
$ws.Cells.Item(409, 17).Value = 18
$ws.Cells.Item(409, 18).Value = "Hortaliza"

# New row 410: "Segunda" quality entry for date 44785 (2022-08-12)
$ws.Cells.Item(410, 1).Value  = 8
$ws.Cells.Item(410, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(410, 3).Value  = "Coquimbo"
$ws.Cells.Item(410, 4).Value  = 44785
$ws.Cells.Item(410, 5).Value  = 4
$ws.Cells.Item(410, 6).Value  = 100112043
$ws.Cells.Item(410, 7).Value  = "Pepino dulce"
$ws.Cells.Item(410, 8).Value  = "Cultivar IV Región"
$ws.Cells.Item(410, 9).Value  = "Segunda"
$ws.Cells.Item(410, 10).Value = 300
$ws.Cells.Item(410, 11).Value = 11000
$ws.Cells.Item(410, 12).Value = 11500
$ws.Cells.Item(410, 13).Value = 11250
$ws.Cells.Item(410, 14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(410, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(410, 16).Value = 625
$ws.Cells.Item(410, 17).Value = 18
$ws.Cells.Item(410, 18).Value = "Hortaliza"
